# 草原建设利用.xlsx update
# 1) Fix the 2020年 row's "草原火灾受害面积" (F10) value: 11045.9 -> 11.0459
# 2) Append a new 2021年 row (row 11) with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the mis-scaled F10 value ---
$ws.Range("F10").Value = 11.0459

# --- Bring the row-11 label/format in line with the other year rows by  ---
# --- copying row 10's A:E formatting (bold/boxed year style) down first ---
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- New row 11 ("2021年") data ---
$ws.Range("A11").Value = "2021年"
$ws.Range("F11").Value = 4.19878
$ws.Range("G11").Value = 7919.25346839773
$ws.Range("H11").Value = 3239.8926
$ws.Range("I11").Value = 37618.8678246127
$ws.Range("J11").Value = 10191.3094666667
